$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 45; this pushes the existing row 45 (old data)
# down to row 46 and duplicates formatting from the row above.
$ws.Rows.Item(45).Insert()

# Row 45 becomes the new weekly entry (newer date / updated price data).
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44890
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112026
$ws.Cells.Item(45, 7).Value = "Haba"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 180
$ws.Cells.Item(45, 11).Value = 9000
$ws.Cells.Item(45, 12).Value = 10000
$ws.Cells.Item(45, 13).Value = 9556
$ws.Cells.Item(45, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región del Maule"
$ws.Cells.Item(45, 16).Value = 382
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
